$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D1").Value = "angkatan"
$ws.Range("E1").Value = "Prodi"
$ws.Range("K1").Value = "HP_Ortu"
$ws.Range("L1").Value = "HP_Mahasiswa"
$ws.Range("B1").Value = "nama-mahasiswa"

$ws.Range("B1").Select() | Out-Null
